$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "obs" column (J) records the observable name "AN" for every data row.
# Update it to "ANep" to reflect the renamed observable (AN -> ANep).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    if ($cell.Value() -eq "AN") {
        $cell.Value = "ANep"
    }
}

# Reflect the last user interaction recorded in the saved file: the
# range J2:J13 (the "obs" column data) ends up selected/active.
$ws.Range("J2:J13").Select()
